{"js": "// Load all body paragraphs (with text) so we can locate the anchors we need by content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Derechos que otorga la conciliaci\u00f3n laboral\" heading paragraph - this paragraph\n// survives the edit (collapsed down to a single empty, underlined paragraph that still carries\n// the _GoBack bookmark).\nlet headingIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(\"Derechos que otorga la conciliaci\u00f3n laboral\") !== -1) {\n        headingIdx = i;\n        break;\n    }\n}\nif (headingIdx === -1) {\n    throw new Error(\"Could not find 'Derechos que otorga la conciliaci\u00f3n laboral' paragraph\");\n}\n\n// Locate the last paragraph that needs to disappear: the blank paragraph that immediately\n// follows \"Para la protecci\u00f3n o asistencia de v\u00edctimas de violencia de g\u00e9nero\".\nlet victimsIdx = -1;\nfor (let i = headingIdx + 1; i < items.length; i++) {\n    if (items[i].text.indexOf(\"Para la protecci\u00f3n o asistencia de v\u00edctimas de violencia de g\u00e9nero\") !== -1) {\n        victimsIdx = i;\n        break;\n    }\n}\nif (victimsIdx === -1) {\n    throw new Error(\"Could not find 'Para la protecci\u00f3n...' paragraph\");\n}\nconst lastRemovedIdx = victimsIdx + 1; // the empty paragraph right after it\n\nconst headingPara = items[headingIdx];\n\n// Apply the surviving paragraph's new formatting (bold off, single underline) while the run\n// with the heading text is still present, so both the run and the paragraph mark pick it up.\nheadingPara.font.bold = false;\nheadingPara.font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark (previously anchored in the paragraph we are about to delete)\n// on this surviving paragraph, scoped to its content so start/end stay inside this paragraph.\nconst contentRange = headingPara.getRange(\"Content\");\ncontentRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Remove the heading text itself, leaving an empty paragraph (but keep the paragraph mark).\nconst found = headingPara.search(\"Derechos que otorga la conciliaci\u00f3n laboral\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n    found.items[0].delete();\n    await context.sync();\n}\n\n// Delete every paragraph between the (now empty) heading paragraph and the trailing blank\n// paragraph, inclusive - walk from the end backwards so indices stay valid.\nparagraphs.load(\"items\");\nawait context.sync();\nfor (let i = lastRemovedIdx; i > headingIdx; i--) {\n    paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Elsewhere in the document, the \"Conclusi\u00f3n y reflexi\u00f3n\" run had a stale\n// lastRenderedPageBreak marker; re-writing the paragraph text regenerates the run cleanly\n// without that artifact while keeping the same visible text/formatting.\nparagraphs.load(\"items/text\");\nawait context.sync();\nfor (const p of paragraphs.items) {\n    if (p.text.indexOf(\"Conclusi\u00f3n y reflexi\u00f3n\") !== -1) {\n        p.insertText(\"Conclusi\u00f3n y reflexi\u00f3n\", Word.InsertLocation.replace);\n        break;\n    }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Derechos que otorga la conciliaci\u00f3n laboral\" heading paragraph - it survives the\n# edit (collapsed down to a single empty, underlined paragraph that still carries the _GoBack\n# bookmark).\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Derechos que otorga la conciliaci\u00f3n laboral*\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not find 'Derechos que otorga la conciliaci\u00f3n laboral' paragraph\"\n}\n\n# Locate the blank paragraph that immediately follows \"Para la protecci\u00f3n o asistencia de\n# v\u00edctimas de violencia de g\u00e9nero\" - everything from the heading through this paragraph\n# (inclusive) collapses away.\n$victimsIndex = -1\nfor ($i = $headingIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Para la protecci\u00f3n o asistencia de v\u00edctimas de violencia de g\u00e9nero*\") {\n        $victimsIndex = $i\n        break\n    }\n}\nif ($victimsIndex -eq -1) {\n    throw \"Could not find 'Para la protecci\u00f3n...' paragraph\"\n}\n$lastRemovedIndex = $victimsIndex + 1\n\n$headingPara = $d.Paragraphs.Item($headingIndex)\n\n# Apply the surviving paragraph's new formatting (bold off, single underline) while the run\n# with the heading text is still present, so both the run and the paragraph mark pick it up.\n$headingPara.Range.Font.Bold = 0\n$headingPara.Range.Font.Underline = 1\n\n# Re-create the \"_GoBack\" bookmark (previously anchored in the paragraph we are about to\n# delete) on this surviving paragraph. Shrink the range by one character first so the bookmark\n# stays inside this paragraph instead of spilling over the paragraph mark.\n$bookmarkRange = $headingPara.Range.Duplicate\n$bookmarkRange.MoveEnd(1, -1) | Out-Null\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n\n# Remove the heading text itself, leaving an empty paragraph (but keep the paragraph mark).\n$textRange = $headingPara.Range.Duplicate\n$textRange.MoveEnd(1, -1) | Out-Null\n$textRange.Text = \"\"\n\n# Delete every paragraph between the (now empty) heading paragraph and the trailing blank\n# paragraph, inclusive - walk from the end backwards so indices stay valid.\nfor ($i = $lastRemovedIndex; $i -gt $headingIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Elsewhere in the document, the \"Conclusi\u00f3n y reflexi\u00f3n\" run had a stale\n# lastRenderedPageBreak marker; re-writing the paragraph text regenerates the run cleanly\n# without that artifact while keeping the same visible text/formatting.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Conclusi\u00f3n y reflexi\u00f3n*\") {\n        $p.Range.Text = \"Conclusi\u00f3n y reflexi\u00f3n\"\n        break\n    }\n}\n\n\"done\"\n"}
